# Update cryptos list figures (prices & 1h volume %) as scraped Fri May 17 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textForceCells = @("D6", "D7", "D12", "D14", "D18", "D20", "D21", "D22", "D23", "D24", "D26", "D29", "D30", "D33", "D35", "D36", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "66.858.61"
$ws.Range("D3").Value = "3.088.88"
$ws.Range("D6").Value = "168.52"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").Value = "3.083.72"
$ws.Range("D12").Value = "0.482"
$ws.Range("D14").Value = "36.36"
$ws.Range("D16").Value = "3.596.96"
$ws.Range("D17").Value = "66.843.57"
$ws.Range("D18").Value = "7.19"
$ws.Range("D19").Value = "3.090.28"
$ws.Range("D20").Value = "16.16"
$ws.Range("D21").Value = "466.33"
$ws.Range("D22").Value = "0.714"
$ws.Range("D23").Value = "7.52"
$ws.Range("D24").Value = "83.74"
$ws.Range("D26").Value = "13.06"
$ws.Range("D29").Value = "8.01"
$ws.Range("D30").Value = "2.39"
$ws.Range("D33").Value = "28.24"
$ws.Range("D35").Value = "1.00"
$ws.Range("D36").Value = "1.01"
$ws.Range("D38").Value = "47.04"
$ws.Range("D40").Value = "50.31"
$ws.Range("D41").Value = "0.318"
$ws.Range("D43").Value = "8.67"
$ws.Range("D44").Value = "2.82"
$ws.Range("D45").Value = "0.0360"
$ws.Range("D46").Value = "382.96"
$ws.Range("D47").Value = "2.778.15"
$ws.Range("D48").Value = "135.12"
$ws.Range("D50").Value = "24.79"
$ws.Range("D51").Value = "2.22"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("E3").Value = "  +5.11%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("E6").Value = "  +5.80%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("E14").Value = "  +6.10%  "
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("E16").Value = "  +4.99%  "
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("E18").Value = "  +4.05%  "
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("E20").Value = "  +7.89%  "
$ws.Range("E21").Value = "  +4.83%  "
$ws.Range("E22").Value = "  +4.06%  "
$ws.Range("E23").Value = "  +3.67%  "
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("E25").Value = "  +6.39%  "
$ws.Range("E26").Value = "  +7.85%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  +3.44%  "
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("E38").Value = "  +5.03%  "
$ws.Range("E39").Value = "  +6.77%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("E50").Value = "  +6.27%  "
$ws.Range("E51").Value = "  +1.86%  "

# --- Rows 40/41 Coin/Link swap ---
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
